$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain literal text
# (matching the source inlineStr cells). Prefixing with an apostrophe forces
# Excel to treat the entry as text instead of auto-converting it to a number,
# and resetting the style back to Normal avoids leaving a stray NumberFormat.

$ws.Range("D2").Value = "'38.335.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.85%  '

$ws.Range("D3").Value = "'2.080.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.55%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'227.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '

$ws.Range("D6").Value = "'0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.75%  '

$ws.Range("D7").Value = "'60.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.69%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +2.20%  '

$ws.Range("D10").Value = "'0.0833"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.68%  '

$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("E12").Value = '  +2.30%  '

$ws.Range("E13").Value = '  +2.68%  '

$ws.Range("D14").Value = "'22.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.89%  '

$ws.Range("D15").Value = "'0.781"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("E16").Value = '  +3.51%  '

$ws.Range("D17").Value = "'2.084.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.27%  '

$ws.Range("D18").Value = "'38.324.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("D19").Value = "'71.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.33%  '

$ws.Range("D20").Value = "'6.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").Value = "'0.0₃0831"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.76%  '

$ws.Range("D22").Value = "'225.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").Value = "'2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.22%  '

$ws.Range("D26").Value = "'169.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").Value = "'9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("D28").Value = "'0.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.05%  '

$ws.Range("D29").Value = "'19.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.77%  '

$ws.Range("E30").Value = '  +8.54%  '

$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").Value = "'2.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.09%  '

$ws.Range("E33").Value = '  +8.06%  '

$ws.Range("D34").Value = "'4.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("E35").Value = '  +0.56%  '

$ws.Range("D36").Value = "'2.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.99%  '

$ws.Range("D37").Value = "'6.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.08%  '

$ws.Range("D38").Value = "'3.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.49%  '

$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.21%  '

$ws.Range("D40").Value = "'18.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.64%  '

$ws.Range("D41").Value = "'1.539.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.10%  '

$ws.Range("D42").Value = "'100.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.61%  '

$ws.Range("E43").Value = '  +2.31%  '

$ws.Range("D44").Value = "'0.0921"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.79%  '

$ws.Range("D45").Value = "'2.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.40%  '

$ws.Range("D46").Value = "'7.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.29%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = "'1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = "'4.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.19%  '

$ws.Range("E49").Value = '  +2.58%  '

$ws.Range("D50").Value = "'2.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.30%  '

$ws.Range("D51").Value = "'2.277.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.38%  '
